$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dates = @("01-10-2021", "02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021")

$r = $lastRow
foreach ($d in $dates) {
    $r = $r + 1
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = "'" + $d
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = 696
    $ws.Cells.Item($r, 3).Value = 3962
    $ws.Cells.Item($r, 4).Value = 59
}
